# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Re-sorts the worker/period detail table (rows 16-50 on Hoja1) so entries
# are ordered by "Periodo Mora" (2107..2111) instead of by worker, folding in
# the first batch of new account-statement periods for each worker.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New values for C (N° Doc Trabajador), D (Nombre Trabajador), E (Periodo Mora),
# F (Valor Mora) and G (Salario Basico) for each row of the table.
$rows = @{
    16 = @("1090462318", "JOSEPH MARQUEZ MARQUEZ VARGAS", "2107", 36341, 908526)
    17 = @("1193561465", "HERNAN DARIO VILA NORIEGA", "2107", 72682, 1817052)
    18 = @("1090985126", "LEONEL CAMILO ESTRADA MARTINEZ", "2107", 81767, 2044184)
    19 = @("45524119", "DANILSA NAVARRO CUETO", "2108", 58146, 1453642)
    20 = @("1052079170", "IDALIDES MARGARITA PONCE RUIZ", "2108", 58146, 1453642)
    21 = @("9282469", "RUBEN SOTO MARTINEZ", "2108", 36341, 908526)
    22 = @("1006887226", "ANGIE VANGRIEKEN JUSAYU", "2108", 42640, 1453642)
    23 = @("73183791", "YESID QUINTANA TORRES", "2108", 36341, 908526)
    24 = @("1090462318", "JOSEPH MARQUEZ MARQUEZ VARGAS", "2108", 36341, 908526)
    25 = @("1193561465", "HERNAN DARIO VILA NORIEGA", "2108", 72682, 1817052)
    26 = @("1090985126", "LEONEL CAMILO ESTRADA MARTINEZ", "2108", 81767, 2044184)
    27 = @("45524119", "DANILSA NAVARRO CUETO", "2109", 58146, 1453642)
    28 = @("1052079170", "IDALIDES MARGARITA PONCE RUIZ", "2109", 58146, 1453642)
    29 = @("9282469", "RUBEN SOTO MARTINEZ", "2109", 36341, 908526)
    30 = @("1006887226", "ANGIE VANGRIEKEN JUSAYU", "2109", 58146, 1453642)
    31 = @("73183791", "YESID QUINTANA TORRES", "2109", 36341, 908526)
    32 = @("1090462318", "JOSEPH MARQUEZ MARQUEZ VARGAS", "2109", 36341, 908526)
    33 = @("1193561465", "HERNAN DARIO VILA NORIEGA", "2109", 72682, 1817052)
    34 = @("1090985126", "LEONEL CAMILO ESTRADA MARTINEZ", "2109", 81767, 2044184)
    35 = @("45524119", "DANILSA NAVARRO CUETO", "2110", 58146, 1453642)
    36 = @("1052079170", "IDALIDES MARGARITA PONCE RUIZ", "2110", 58146, 1453642)
    37 = @("9282469", "RUBEN SOTO MARTINEZ", "2110", 36341, 908526)
    38 = @("1006887226", "ANGIE VANGRIEKEN JUSAYU", "2110", 58146, 1453642)
    39 = @("73183791", "YESID QUINTANA TORRES", "2110", 36341, 908526)
    40 = @("1090462318", "JOSEPH MARQUEZ MARQUEZ VARGAS", "2110", 36341, 908526)
    41 = @("1193561465", "HERNAN DARIO VILA NORIEGA", "2110", 72682, 1817052)
    42 = @("1090985126", "LEONEL CAMILO ESTRADA MARTINEZ", "2110", 81767, 2044184)
    43 = @("45524119", "DANILSA NAVARRO CUETO", "2111", 44579, 1453642)
    44 = @("1052079170", "IDALIDES MARGARITA PONCE RUIZ", "2111", 44579, 1453642)
    45 = @("9282469", "RUBEN SOTO MARTINEZ", "2111", 27861, 908526)
    46 = @("1006887226", "ANGIE VANGRIEKEN JUSAYU", "2111", 44579, 1453642)
    47 = @("73183791", "YESID QUINTANA TORRES", "2111", 27861, 908526)
    48 = @("1090462318", "JOSEPH MARQUEZ MARQUEZ VARGAS", "2111", 27861, 908526)
    49 = @("1193561465", "HERNAN DARIO VILA NORIEGA", "2111", 55723, 1817052)
    50 = @("1090985126", "LEONEL CAMILO ESTRADA MARTINEZ", "2111", 62688, 2044184)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("C$r").Value = $vals[0]
    $ws.Range("D$r").Value = $vals[1]
    $ws.Range("E$r").Value = $vals[2]
    $ws.Range("F$r").Value = $vals[3]
    $ws.Range("G$r").Value = $vals[4]
}
